$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-detected as a Number by Excel;
# force them to stay Text (matching original inlineStr type) without leaving a
# lingering custom number format on the cell.
$textCells = @('D5', 'D6', 'D7', 'D8', 'D10', 'D11', 'D13', 'D15', 'D16', 'D18', 'D19', 'D22', 'D23', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D33', 'D35', 'D36', 'D38', 'D39', 'D42', 'D43', 'D44', 'D48', 'D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply every cell value change from the diff, row by row.
$ws.Range('D2').Value = '30.518.57'
$ws.Range('E2').Value = '  +1.97%  '
$ws.Range('D3').Value = '1.673.66'
$ws.Range('E3').Value = '  +2.46%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '219.83'
$ws.Range('E5').Value = '  +2.53%  '
$ws.Range('D6').Value = '0.531'
$ws.Range('E6').Value = '  +2.64%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '29.70'
$ws.Range('E8').Value = '  +4.02%  '
$ws.Range('E9').Value = '  +2.82%  '
$ws.Range('D10').Value = '0.0641'
$ws.Range('E10').Value = '  +5.52%  '
$ws.Range('D11').Value = '0.0906'
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('D12').Value = '1.913.88'
$ws.Range('E12').Value = '  +2.44%  '
$ws.Range('D13').Value = '0.615'
$ws.Range('E13').Value = '  +9.25%  '
$ws.Range('D14').Value = '1.665.46'
$ws.Range('D15').Value = '10.20'
$ws.Range('E15').Value = '  +9.13%  '
$ws.Range('D16').Value = '3.99'
$ws.Range('E16').Value = '  +3.59%  '
$ws.Range('D17').Value = '30.537.89'
$ws.Range('E17').Value = '  +1.97%  '
$ws.Range('D18').Value = '66.50'
$ws.Range('E18').Value = '  +3.78%  '
$ws.Range('D19').Value = '242.91'
$ws.Range('E19').Value = '  +0.28%  '
$ws.Range('D20').Value = '0.0₃0723'
$ws.Range('E20').Value = '  +3.17%  '
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('D22').Value = '4.28'
$ws.Range('E22').Value = '  +3.48%  '
$ws.Range('D23').Value = '10.00'
$ws.Range('E23').Value = '  +1.70%  '
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').Value = '158.51'
$ws.Range('E25').Value = '  +0.49%  '
$ws.Range('D26').Value = '15.86'
$ws.Range('E26').Value = '  +2.20%  '
$ws.Range('D27').Value = '0.113'
$ws.Range('E27').Value = '  +2.43%  '
$ws.Range('D28').Value = '6.68'
$ws.Range('E28').Value = '  +0.98%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('D30').Value = '0.0497'
$ws.Range('E30').Value = '  +2.05%  '
$ws.Range('E31').Value = '  +2.96%  '
$ws.Range('E32').Value = '  +2.82%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '3.28'
$ws.Range('E33').Value = '  +3.36%  '
$ws.Range('B34').Value = 'Maker'
$ws.Range('C34').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D34').Value = '1.497.56'
$ws.Range('E34').Value = '  +5.07%  '
$ws.Range('D35').Value = '1.77'
$ws.Range('E35').Value = '  +7.11%  '
$ws.Range('D36').Value = '84.25'
$ws.Range('E36').Value = '  +10.74%  '
$ws.Range('E37').Value = '  -0.79%  '
$ws.Range('D38').Value = '0.599'
$ws.Range('E38').Value = '  +8.25%  '
$ws.Range('D39').Value = '0.0178'
$ws.Range('E39').Value = '  +5.02%  '
$ws.Range('E40').Value = '  -4.46%  '
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('D42').Value = '0.839'
$ws.Range('E42').Value = '  +1.08%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').Value = '0.0499'
$ws.Range('E43').Value = '  +1.72%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '1.97'
$ws.Range('E44').Value = '  -1.47%  '
$ws.Range('E45').Value = '  +0.30%  '
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('E47').Value = '  +3.38%  '
$ws.Range('D48').Value = '50.98'
$ws.Range('E48').Value = '  -3.73%  '
$ws.Range('D49').Value = '1.806.00'
$ws.Range('E49').Value = '  +1.74%  '
$ws.Range('D50').Value = '94.86'
$ws.Range('E50').Value = '  +4.73%  '
$ws.Range('D51').Value = '0.0₆0114'
$ws.Range('E51').Value = '  +0.91%  '

# Restore the default 'Normal' style on the forced-text cells so no stray
# number-format style index is left behind.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
